# Update the LinkedIn export table with the new data set.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Franck Bonnier"
$ws.Range("B2").Value = "Head of IT Transverse Department chez ArcelorMittal France"
$ws.Range("C2").Value = "ArcelorMittal"
$ws.Range("D2").Value = "France"
$ws.Range("E2").Value = "https://www.linkedin.com/in/franck-bonnier-04033221"
$ws.Range("F2").Value = "https://www.linkedin.com/in/hamouda-makhloufi-934a7437"
$ws.Range("A3").Value = "Denis Dey"
$ws.Range("B3").Value = "Responsable de programme data"
$ws.Range("C3").Value = "ArcelorMittal"
$ws.Range("D3").Value = "Lille, Hauts-de-France, France"
$ws.Range("E3").Value = "https://www.linkedin.com/in/denis-dey-88167528"
$ws.Range("F3").Value = "https://www.linkedin.com/in/hamouda-makhloufi-934a7437"
$ws.Range("A4").Value = "Chankar Pourouchotamane"
$ws.Range("B4").Value = "Head of Data Office Group"
$ws.Range("C4").Value = "CMA CGM"
$ws.Range("D4").Value = "Greater Marseille Metropolitan Area"
$ws.Range("E4").Value = "https://www.linkedin.com/in/chankar-pourouchotamane-79a5687"
$ws.Range("F4").Value = "https://www.linkedin.com/in/hamouda-makhloufi-934a7437"
$ws.Range("A5").Value = "Gabriel Fricout"
$ws.Range("B5").Value = "Head Of Data Intelligence Department chez ArcelorMittal France"
$ws.Range("C5").Value = "ArcelorMittal France"
$ws.Range("D5").Value = "Greater Metz Area"
$ws.Range("E5").Value = "https://www.linkedin.com/in/gabriel-fricout-35709913"
$ws.Range("F5").Value = "https://www.linkedin.com/in/hamouda-makhloufi-934a7437"
$ws.Range("A6").Value = "Tharaud J."
$ws.Range("B6").Value = "Chief Data Officer / Head of Data Plateform Engineering"
$ws.Range("C6").Value = "PMU"
$ws.Range("D6").Value = "Brunoy, Île-de-France, France"
$ws.Range("E6").Value = "https://www.linkedin.com/in/jtharaud"
$ws.Range("F6").Value = "https://www.linkedin.com/in/hamouda-makhloufi-934a7437"
$ws.Range("A7").Value = "Matthieu Bonan"
$ws.Range("B7").Value = "Head of Data Office chez Groupe KILOUTOU"
$ws.Range("C7").Value = "Groupe KILOUTOU"
$ws.Range("D7").Value = "Lille, Hauts-de-France, France"
$ws.Range("E7").Value = "https://www.linkedin.com/in/matthieu-bonan-aa8263109"
$ws.Range("F7").Value = "https://www.linkedin.com/in/hamouda-makhloufi-934a7437"
$ws.Range("A8").Value = "Naceur Abderrahim"
$ws.Range("B8").Value = "Head of Data Management and Governance"
$ws.Range("C8").Value = "Servier"
$ws.Range("D8").Value = "Paris, Île-de-France, France"
$ws.Range("E8").Value = "https://www.linkedin.com/in/naceur-abderrahim-25b6b29"
$ws.Range("F8").Value = "https://www.linkedin.com/in/hamouda-makhloufi-934a7437"
$ws.Range("A9").Value = "Thierry Mocquillon"
$ws.Range("B9").Value = "Directeur Technologie et Système d’Information - Infrastructures, Technologies et Services"
$ws.Range("C9").Value = "Groupe Covéa"
$ws.Range("D9").Value = "La Rochelle, Nouvelle-Aquitaine, France"
$ws.Range("E9").Value = "https://www.linkedin.com/in/thierry-mocquillon-39636236"
$ws.Range("F9").Value = "https://www.linkedin.com/in/alexandre-rouger-216899198"
$ws.Range("A10").Value = "Cyprien Falque"
$ws.Range("B10").Value = "Directeur général / CEO chez S3NS"
$ws.Range("C10").Value = "S3NS"
$ws.Range("D10").Value = "Greater Paris Metropolitan Region"
$ws.Range("E10").Value = "https://www.linkedin.com/in/cyprien-falque-05878131"
$ws.Range("F10").Value = "https://www.linkedin.com/in/olivier-esposito-861464155"
$ws.Range("A11").Value = "Mickael DUPONT"
$ws.Range("B11").Value = "Actuaire / Data Scientist / Manager IT"
$ws.Range("C11").Value = "Groupe Covéa"
$ws.Range("D11").Value = "Niort, Nouvelle-Aquitaine, France"
$ws.Range("E11").Value = "https://www.linkedin.com/in/mickael-dupont-61950810"
$ws.Range("F11").Value = "https://www.linkedin.com/in/david-brenet-ab26367b"
$ws.Range("A12").Value = "Mehdy A."
$ws.Range("B12").Value = "Responsable Pôle Infrastructure chez APIVIA Courtage, filiale d'Apivia Macif Mutuelle - Aéma Groupe"
$ws.Range("C12").Value = "Apivia Courtage"
$ws.Range("D12").Value = "Poitiers, Nouvelle-Aquitaine, France"
$ws.Range("E12").Value = "https://www.linkedin.com/in/mehdy-a-4a865baa"
$ws.Range("F12").Value = "https://www.linkedin.com/in/cedric-porchier-43ba2a194"
$ws.Range("A13").Value = "emmanuel Roquigny"
$ws.Range("B13").Value = "Responsable SIRH"
$ws.Range("C13").Value = "Inter Mutuelle Assistance"
$ws.Range("D13").Value = "Greater La Rochelle Area"
$ws.Range("E13").Value = "https://www.linkedin.com/in/emmanuel-roquigny-217521b6"
$ws.Range("F13").Value = "https://www.linkedin.com/in/fabien-leroy-70683070"
$ws.Range("A14").Value = "Olivier Louis MONNIER"
$ws.Range("B14").Value = "Group Chief Data & AI Officer"
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = "Paris, Île-de-France, France"
$ws.Range("E14").Value = "https://www.linkedin.com/in/olivier-louis-monnier"
$ws.Range("F14").Value = "https://www.linkedin.com/in/leverageondata"
$ws.Range("A15").Value = "Arnaud HAMEL"
$ws.Range("B15").Value = "Head of IT Architecture chez Rexel"
$ws.Range("C15").Value = "Rexel"
$ws.Range("D15").Value = "Greater Paris Metropolitan Region"
$ws.Range("E15").Value = "https://www.linkedin.com/in/arnaud-hamel-a1902665/"
$ws.Range("F15").Value = ""
$ws.Range("A16").Value = "Jean Dupont (Mock)"
$ws.Range("B16").Value = "Développeur Python Senior"
$ws.Range("C16").Value = "Mock Corp"
$ws.Range("D16").Value = "Paris, France"
$ws.Range("E16").Value = "https://www.linkedin.com/in/nicolas-d-avout-d-auerstaedt-16201b53"
$ws.Range("F16").Value = "https://www.linkedin.com/in/fabriceleyglene"
$ws.Range("A17").Value = "Jean Dupont (Mock)"
$ws.Range("B17").Value = "Développeur Python Senior"
$ws.Range("C17").Value = "Mock Corp"
$ws.Range("D17").Value = "Paris, France"
$ws.Range("E17").Value = "https://www.linkedin.com/in/jvaucelle"
$ws.Range("F17").Value = "https://www.linkedin.com/in/laetitia-palatini-a3b49412b"
$ws.Range("A18").Value = "Jean Dupont (Mock)"
$ws.Range("B18").Value = "Développeur Python Senior"
$ws.Range("C18").Value = "Mock Corp"
$ws.Range("D18").Value = "Paris, France"
$ws.Range("E18").Value = "https://www.linkedin.com/in/jacques-chemaoun-43253790"
$ws.Range("F18").Value = "https://www.linkedin.com/in/nicolas-d-avout-d-auerstaedt-16201b53"
$ws.Range("A19").Value = "Jean Dupont (Mock)"
$ws.Range("B19").Value = "Développeur Python Senior"
$ws.Range("C19").Value = "Mock Corp"
$ws.Range("D19").Value = "Paris, France"
$ws.Range("E19").Value = "https://www.linkedin.com/in/sandrine-racouchot-74aa0172"
$ws.Range("F19").Value = "https://www.linkedin.com/in/nicolas-d-avout-d-auerstaedt-16201b53"
$ws.Range("A20").Value = "Jean Dupont (Mock)"
$ws.Range("B20").Value = "Développeur Python Senior"
$ws.Range("C20").Value = "Mock Corp"
$ws.Range("D20").Value = "Paris, France"
$ws.Range("E20").Value = "https://www.linkedin.com/in/cyrilbayon"
$ws.Range("F20").Value = "https://www.linkedin.com/in/nicolas-d-avout-d-auerstaedt-16201b53"

Write-Host "Updated rows 2-20 of sheet $($ws.Name)"
